$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "289.23"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.42%"
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.24"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-7.18%"
$ws.Range("E3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.912"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.83%"
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07243"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.12%"
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.787"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-14.26%"
$ws.Range("E6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.576"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.95%"
$ws.Range("E7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.709"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.24%"
$ws.Range("E8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9013"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.40%"
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1672"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.24%"
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07823"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.97%"
$ws.Range("E11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07982"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.27%"
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03038"
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.79%"
$ws.Range("E13").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.34%"
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001496"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.65%"
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005758"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.21%"
$ws.Range("E16").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.01%"
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.074"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.85%"
$ws.Range("E18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3320"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.09%"
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1295"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.37%"
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.963"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-9.61%"
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2176"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.18%"
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04499"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.89%"
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001213"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.62%"
$ws.Range("E24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004426"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.69%"
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("D26").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.70%"
$ws.Range("E26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003390"
$ws.Range("D27").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-95.48%"
$ws.Range("E27").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01569"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-9.14%"
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04333"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.91%"
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007328"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.60%"
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01006"
$ws.Range("D42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1304"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.27%"
$ws.Range("E43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002003"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.19%"
$ws.Range("E44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009031"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-14.94%"
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005894"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.87%"
$ws.Range("E46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("E47").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "166.52%"
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002895"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.49%"
$ws.Range("E49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("E50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.11%"
$ws.Range("E51").ClearFormats()
